$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (F column) values for several rows.
# (values re-pulled / recalculated per commit message "repull data, push all data, mean calculation")
$ws.Range("F2").Value  = -1
$ws.Range("F3").Value  = -6
$ws.Range("F5").Value  = -2
$ws.Range("F9").Value  = -7
$ws.Range("F10").Value = -5
$ws.Range("F13").Value = -4
$ws.Range("F19").Value = -4
$ws.Range("F26").Value = 0
$ws.Range("F31").Value = -2
$ws.Range("F33").Value = -2
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F41").Value = -1
$ws.Range("F45").Value = 1
$ws.Range("F46").Value = 3
